$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change C2's text to "Install".
$ws.Range("C2").Value = "Install"

# Add an internal-style hyperlink on C2 pointing at the itms-services
# install manifest (no external relationship - this goes through
# SubAddress, mirroring the `location=` attribute in the OOXML).
$ws.Hyperlinks.Add(
    $ws.Range("C2"),
    "",
    "itms-services://?action=download-manifest&url=https://a-randomm-user.github.io/IPA-save/onlineDownload/official-minecraftpe-479516143/0.1.2.plist",
    "",
    "Install"
)

# Make C2 look like the other hyperlink cells (D2/E2/F2 already use the
# "Hyperlink" cell style) by copying their formatting rather than
# re-applying a named style (which would create a brand-new cell format).
$ws.Range("D2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column C no longer needs to be wide enough to show the raw URL.
$ws.Columns.Item(3).ColumnWidth = 6.1875

# Move the active selection.
$ws.Range("D6").Select()
